# Second checkpoint commit - Saturday session
# Applies the "week 10" task-table updates:
#  - mark "sitemap updated" (row 21) as done ("y")
#  - insert a new task row for "javascript buttons to hide", done, assigned to Derek
#  - move the selection / view to reflect the edit location

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("task progress")

# Mark row 21 ("sitemap updated") as done.
$ws.Range("C21").Value = "y"

# Row 25 was an empty gap row between 24 ("procedural navbar...") and 26
# ("Let Emer know..."); fill it in with the new task - no row insertion/shift.
$ws.Range("B25").Value = "javascript buttons to hide"
$ws.Range("C25").Value = "y"
$ws.Range("D25").Value = "Derek"

# Columns B, D, E were manually narrowed (no longer auto "best fit"); apply
# the new fixed widths (closest achievable via the character-unit ColumnWidth
# property to the saved OOXML widths 107.85546875 / 120.85546875 / 46.28515625).
$ws.Columns.Item(2).ColumnWidth = 107
$ws.Columns.Item(4).ColumnWidth = 120
$ws.Columns.Item(5).ColumnWidth = 45.5

# Update the view: scroll/selection moved to C21.
$ws.Range("C21").Select()

$wb.Save()
